$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Progress Tracker")

# Week 8 lectures completed (rows 46-50): mark status COMPLETE and record
# the completion date, reusing the existing date-cell formatting (style
# of B44) so no new styles are introduced. Row 51 moves to IN PROGRESS.

$ws.Range("B44").Copy()

$ws.Range("A46").Value = "COMPLETE"
$ws.Range("B46").PasteSpecial(-4122)
$ws.Range("B46").Value = "1/10/2021"

$ws.Range("A47").Value = "COMPLETE"
$ws.Range("B47").PasteSpecial(-4122)
$ws.Range("B47").Value = "1/10/2021"

$ws.Range("A48").Value = "COMPLETE"
$ws.Range("B48").PasteSpecial(-4122)
$ws.Range("B48").Value = "1/10/2021"

$ws.Range("A49").Value = "COMPLETE"
$ws.Range("B49").PasteSpecial(-4122)
$ws.Range("B49").Value = "1/10/2021"

$ws.Range("A50").Value = "COMPLETE"
$ws.Range("B50").PasteSpecial(-4122)
$ws.Range("B50").Value = "1/13/2021"

$ws.Range("A51").Value = "IN PROGRESS"

# Update the active selection to follow where the user left off.
$ws.Range("B51").Select()
